# Update docx golden tests for reference doc changes.
#
# 1. Remove Font.Spacing / Font.Kerning overrides from the "Title" and
#    "TitleChar" styles (they revert to their base style defaults).
# 2. Make "Author" and "Date" styles based on "Title" (instead of
#    implicitly on Normal), drop the explicit center alignment (now it
#    comes from the Title base style), and give them their own run
#    size of 12pt (sz/szCs = 24 half-points) instead of Title's 28pt.

$d = $word.ActiveDocument

# --- Title / TitleChar: drop the -10 character spacing and 28pt kerning ---
foreach ($styleName in @("Title", "TitleChar")) {
    $s = $d.Styles($styleName)
    $s.Font.Spacing = 0
    $s.Font.Kerning = 0
}

# --- Author style: base on Title, remove explicit centering, set 12pt run size ---
$author = $d.Styles("Author")
$author.BaseStyle = "Title"
$author.ParagraphFormat.Alignment = 0
$author.Font.Size = 12

# --- Date style: base on Title, remove explicit centering, set 12pt run size ---
$date = $d.Styles("Date")
$date.BaseStyle = "Title"
$date.ParagraphFormat.Alignment = 0
$date.Font.Size = 12
